$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire "publisher_id" column (column E) — its values are no
# longer needed now that publishers are looked up via a hashmap instead.
$ws.Range("E1").EntireColumn.Delete()

# Put the active selection on E1 (first cell of what is now "publisher_name").
$ws.Range("E1").Select()
